# Fixing salinity correction following Weiss 1974 and Wiesenburg and Guinasso 1979
# Updates recalculated ch4_mmolm3, d13_ch4_permil, co2_mmolm3, d13_co2_permil values
# for rows 2-7 (columns Y, Z, AA, AB).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "Y2"  = 0.02770131288012934
    "Z2"  = -75.44098869796578
    "AA2" = 21.06555855794637
    "AB2" = -11.35871179134607

    "Y3"  = 0.06275479364418075
    "Z3"  = -71.0924594017186
    "AA3" = 60.93540747940875
    "AB3" = -14.40110551630081

    "Y4"  = 0.06827144592748455
    "Z4"  = -70.47150379399436
    "AA4" = 62.84622280740862
    "AB4" = -13.94303078788584

    "Y5"  = 0.07756118444500827
    "Z5"  = -71.8121046011868
    "AA5" = 63.98566571139974
    "AB5" = -14.47375611850876

    "Y6"  = 0.08025936103358021
    "Z6"  = -71.61804818071116
    "AA6" = 63.87363607013693
    "AB6" = -14.06353980954156

    "Y7"  = 0.06856818159211989
    "Z7"  = -72.02592110609297
    "AA7" = 65.83034680576739
    "AB7" = -14.39355243823259
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
